# Auto-generated Excel COM-interop script
# Applies numeric corrections to Leve profit-calculation columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Row=87; Col=8; Value=23778.5}
    @{Row=87; Col=10; Value=23778.5}
    @{Row=87; Col=12; Value=23778.5}
    @{Row=87; Col=14; Value=-26274.5}
    @{Row=90; Col=8; Value=23778.5}
    @{Row=90; Col=10; Value=23778.5}
    @{Row=90; Col=12; Value=71335.5}
    @{Row=90; Col=14; Value=-83815.5}
    @{Row=98; Col=8; Value=4784.3887}
    @{Row=98; Col=9; Value=3226.5833}
    @{Row=98; Col=10; Value=7900}
    @{Row=98; Col=11; Value=3226.5833}
    @{Row=98; Col=12; Value=7900}
    @{Row=98; Col=13; Value=-1728.5833}
    @{Row=98; Col=14; Value=-10896}
    @{Row=122; Col=8; Value=4784.3887}
    @{Row=122; Col=9; Value=3226.5833}
    @{Row=122; Col=10; Value=7900}
    @{Row=122; Col=11; Value=9679.749899999999}
    @{Row=122; Col=12; Value=23700}
    @{Row=122; Col=13; Value=-7229.749899999999}
    @{Row=122; Col=14; Value=-28600}
    @{Row=137; Col=8; Value=3423.1226}
    @{Row=137; Col=9; Value=2787.5366}
    @{Row=137; Col=10; Value=6680.5}
    @{Row=137; Col=11; Value=8362.6098}
    @{Row=137; Col=12; Value=20041.5}
    @{Row=137; Col=13; Value=-5812.6098}
    @{Row=137; Col=14; Value=-25141.5}
    @{Row=138; Col=8; Value=2312.6}
    @{Row=138; Col=10; Value=2949.3906}
    @{Row=138; Col=12; Value=8848.1718}
    @{Row=138; Col=14; Value=-19128.1718}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Row=2; Col=8; Value=2584.1428}
    @{Row=2; Col=9; Value=1372.75}
    @{Row=2; Col=10; Value=4199.3335}
    @{Row=2; Col=11; Value=1372.75}
    @{Row=2; Col=12; Value=4199.3335}
    @{Row=2; Col=13; Value=-1259.75}
    @{Row=2; Col=14; Value=-4425.3335}
    @{Row=10; Col=8; Value=17191.5}
    @{Row=10; Col=10; Value=17191.5}
    @{Row=10; Col=12; Value=17191.5}
    @{Row=10; Col=14; Value=-17531.5}
    @{Row=32; Col=8; Value=9278.266}
    @{Row=32; Col=9; Value=6310.8525}
    @{Row=32; Col=11; Value=6310.8525}
    @{Row=32; Col=13; Value=-6023.8525}
    @{Row=63; Col=8; Value=7293964.5}
    @{Row=63; Col=9; Value=19789688}
    @{Row=63; Col=10; Value=4792.6665}
    @{Row=63; Col=11; Value=19789688}
    @{Row=63; Col=12; Value=4792.6665}
    @{Row=63; Col=13; Value=-19789002}
    @{Row=63; Col=14; Value=-6164.6665}
    @{Row=66; Col=8; Value=7293964.5}
    @{Row=66; Col=9; Value=19789688}
    @{Row=66; Col=10; Value=4792.6665}
    @{Row=66; Col=11; Value=98948440}
    @{Row=66; Col=12; Value=23963.3325}
    @{Row=66; Col=13; Value=-98945008}
    @{Row=66; Col=14; Value=-30827.3325}
    @{Row=102; Col=9; Value=2400}
    @{Row=102; Col=10; Value=2200}
    @{Row=102; Col=11; Value=2400}
    @{Row=102; Col=12; Value=2200}
    @{Row=102; Col=13; Value=-778}
    @{Row=102; Col=14; Value=-5444}
    @{Row=116; Col=8; Value=2584.1428}
    @{Row=116; Col=9; Value=1372.75}
    @{Row=116; Col=10; Value=4199.3335}
    @{Row=116; Col=11; Value=1372.75}
    @{Row=116; Col=12; Value=4199.3335}
    @{Row=116; Col=13; Value=921.25}
    @{Row=116; Col=14; Value=-8787.333500000001}
    @{Row=122; Col=8; Value=2314.4194}
    @{Row=122; Col=9; Value=1351.6086}
    @{Row=122; Col=10; Value=5082.5}
    @{Row=122; Col=11; Value=4054.8258}
    @{Row=122; Col=12; Value=15247.5}
    @{Row=122; Col=13; Value=-1604.8258}
    @{Row=122; Col=14; Value=-20147.5}
    @{Row=132; Col=8; Value=2427.0312}
    @{Row=132; Col=9; Value=1291.0834}
    @{Row=132; Col=11; Value=3873.2502}
    @{Row=132; Col=13; Value=-1343.2502}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{Row=3; Col=8; Value=2584.1428}
    @{Row=3; Col=9; Value=1372.75}
    @{Row=3; Col=10; Value=4199.3335}
    @{Row=3; Col=11; Value=1372.75}
    @{Row=3; Col=12; Value=4199.3335}
    @{Row=3; Col=13; Value=-1258.75}
    @{Row=3; Col=14; Value=-4427.3335}
    @{Row=138; Col=8; Value=42682.285}
    @{Row=138; Col=10; Value=42682.285}
    @{Row=138; Col=12; Value=42682.285}
    @{Row=138; Col=14; Value=-52962.285}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Row=99; Col=8; Value=8666.666999999999}
    @{Row=99; Col=9; Value=1000}
    @{Row=99; Col=10; Value=12500}
    @{Row=99; Col=11; Value=1000}
    @{Row=99; Col=12; Value=12500}
    @{Row=99; Col=13; Value=498}
    @{Row=99; Col=14; Value=-15496}
    @{Row=126; Col=8; Value=8666.666999999999}
    @{Row=126; Col=9; Value=1000}
    @{Row=126; Col=10; Value=12500}
    @{Row=126; Col=11; Value=3000}
    @{Row=126; Col=12; Value=37500}
    @{Row=126; Col=13; Value=-530}
    @{Row=126; Col=14; Value=-42440}
    @{Row=130; Col=8; Value=43280}
    @{Row=130; Col=10; Value=43280}
    @{Row=130; Col=12; Value=43280}
    @{Row=130; Col=14; Value=-53320}
    @{Row=134; Col=8; Value=6336.68}
    @{Row=134; Col=10; Value=4840}
    @{Row=134; Col=12; Value=14520}
    @{Row=134; Col=14; Value=-19590}
    @{Row=138; Col=8; Value=47823.6}
    @{Row=138; Col=10; Value=47823.6}
    @{Row=138; Col=12; Value=47823.6}
    @{Row=138; Col=14; Value=-58103.6}
    @{Row=140; Col=8; Value=118943.75}
    @{Row=140; Col=10; Value=118943.75}
    @{Row=140; Col=12; Value=118943.75}
    @{Row=140; Col=14; Value=-129303.75}
    @{Row=141; Col=8; Value=15390}
    @{Row=141; Col=10; Value=15390}
    @{Row=141; Col=12; Value=15390}
    @{Row=141; Col=14; Value=-25750}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Row=2; Col=8; Value=1785791.8}
    @{Row=2; Col=9; Value=57.142857}
    @{Row=2; Col=10; Value=2747341.2}
    @{Row=2; Col=11; Value=342.857142}
    @{Row=2; Col=12; Value=16484047.2}
    @{Row=2; Col=13; Value=-229.857142}
    @{Row=2; Col=14; Value=-16484273.2}
    @{Row=33; Col=8; Value=104.14286}
    @{Row=33; Col=9; Value=99.666664}
    @{Row=33; Col=10; Value=107.5}
    @{Row=33; Col=11; Value=597.999984}
    @{Row=33; Col=12; Value=645}
    @{Row=33; Col=13; Value=-314.999984}
    @{Row=33; Col=14; Value=-1211}
    @{Row=38; Col=8; Value=153.375}
    @{Row=38; Col=9; Value=81.666664}
    @{Row=38; Col=10; Value=196.4}
    @{Row=38; Col=11; Value=244.999992}
    @{Row=38; Col=12; Value=589.2}
    @{Row=38; Col=13; Value=102.000008}
    @{Row=38; Col=14; Value=-1283.2}
    @{Row=51; Col=8; Value=2999.875}
    @{Row=51; Col=9; Value=0}
    @{Row=51; Col=10; Value=2999.875}
    @{Row=51; Col=11; Value=0}
    @{Row=51; Col=12; Value=8999.625}
    @{Row=51; Col=14; Value=-9919.625}
    @{Row=55; Col=8; Value=252501.25}
    @{Row=55; Col=10; Value=4002.5}
    @{Row=55; Col=12; Value=12007.5}
    @{Row=55; Col=14; Value=-12361.5}
    @{Row=106; Col=8; Value=4737.5}
    @{Row=106; Col=10; Value=4737.5}
    @{Row=106; Col=12; Value=14212.5}
    @{Row=106; Col=14; Value=-16104.5}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Row=9; Col=8; Value=10266.333}
    @{Row=9; Col=9; Value=799}
    @{Row=9; Col=10; Value=15000}
    @{Row=9; Col=11; Value=799}
    @{Row=9; Col=12; Value=15000}
    @{Row=9; Col=13; Value=-629}
    @{Row=9; Col=14; Value=-15340}
    @{Row=52; Col=8; Value=31332.223}
    @{Row=52; Col=10; Value=31332.223}
    @{Row=52; Col=12; Value=31332.223}
    @{Row=52; Col=14; Value=-31850.223}
    @{Row=102; Col=8; Value=2533.1875}
    @{Row=102; Col=9; Value=2135.4}
    @{Row=102; Col=11; Value=2135.4}
    @{Row=102; Col=13; Value=-513.4000000000001}
    @{Row=122; Col=8; Value=3241.6667}
    @{Row=122; Col=9; Value=2544.7334}
    @{Row=122; Col=10; Value=6726.3335}
    @{Row=122; Col=11; Value=7634.2002}
    @{Row=122; Col=12; Value=20179.0005}
    @{Row=122; Col=13; Value=-5184.2002}
    @{Row=122; Col=14; Value=-25079.0005}
    @{Row=126; Col=8; Value=3637.01}
    @{Row=126; Col=9; Value=2770.0352}
    @{Row=126; Col=10; Value=4786.256}
    @{Row=126; Col=11; Value=8310.105599999999}
    @{Row=126; Col=12; Value=14358.768}
    @{Row=126; Col=13; Value=-5840.105599999999}
    @{Row=126; Col=14; Value=-19298.768}
    @{Row=132; Col=8; Value=2692.7778}
    @{Row=132; Col=9; Value=1171.0625}
    @{Row=132; Col=10; Value=4906.1816}
    @{Row=132; Col=11; Value=3513.1875}
    @{Row=132; Col=12; Value=14718.5448}
    @{Row=132; Col=13; Value=-983.1875}
    @{Row=132; Col=14; Value=-19778.5448}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Row=40; Col=8; Value=5316.081}
    @{Row=40; Col=9; Value=5087}
    @{Row=40; Col=10; Value=6028.778}
    @{Row=40; Col=11; Value=5087}
    @{Row=40; Col=12; Value=6028.778}
    @{Row=40; Col=13; Value=-4951}
    @{Row=40; Col=14; Value=-6300.778}
    @{Row=75; Col=8; Value=46875}
    @{Row=75; Col=10; Value=46875}
    @{Row=75; Col=12; Value=46875}
    @{Row=75; Col=14; Value=-48747}
    @{Row=78; Col=8; Value=46875}
    @{Row=78; Col=10; Value=46875}
    @{Row=78; Col=12; Value=140625}
    @{Row=78; Col=14; Value=-149985}
    @{Row=122; Col=8; Value=5503.96}
    @{Row=122; Col=9; Value=4146.077}
    @{Row=122; Col=11; Value=12438.231}
    @{Row=122; Col=13; Value=-9988.231}
    @{Row=132; Col=8; Value=6056.6855}
    @{Row=132; Col=9; Value=2312.6}
    @{Row=132; Col=10; Value=7554.32}
    @{Row=132; Col=11; Value=6937.799999999999}
    @{Row=132; Col=12; Value=22662.96}
    @{Row=132; Col=13; Value=-4407.799999999999}
    @{Row=132; Col=14; Value=-27722.96}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Row=94; Col=8; Value=35000}
    @{Row=94; Col=10; Value=35000}
    @{Row=94; Col=12; Value=35000}
    @{Row=94; Col=14; Value=-36802}
    @{Row=122; Col=8; Value=6772}
    @{Row=122; Col=9; Value=5437.091}
    @{Row=122; Col=11; Value=16311.273}
    @{Row=122; Col=13; Value=-13861.273}
    @{Row=126; Col=8; Value=720806.4399999999}
    @{Row=126; Col=9; Value=4600}
    @{Row=126; Col=11; Value=13800}
    @{Row=126; Col=13; Value=-11330}
    @{Row=132; Col=8; Value=10106328}
    @{Row=132; Col=9; Value=9306.666999999999}
    @{Row=132; Col=10; Value=15876054}
    @{Row=132; Col=11; Value=27920.001}
    @{Row=132; Col=12; Value=47628162}
    @{Row=132; Col=13; Value=-25390.001}
    @{Row=132; Col=14; Value=-47633222}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# --- Cell removed entirely (no longer present after the edit) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 13).ClearContents()
